# The "Strike#" column (G, header "K") was regenerated from freshly
# recomputed s_vals. Write the new K values for rows 2-24 in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = 3
    15 = 1
    16 = 3
    17 = 1
    18 = 1
    19 = 6
    20 = 1
    21 = 5
    22 = 1
    23 = 2
    24 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
